$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.521.27'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.944.50'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '571.01'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.79'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '2.941.40'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('E11').Value = '  -3.85%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.463'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.82'
$ws.Range('D14').Style = "Normal"
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '65.549.98'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '3.434.81'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').Value = '2.946.24'
$ws.Range('E19').Value = '  -2.01%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.66'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +12.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '445.53'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.28'
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.09'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('E25').Value = '  -1.33%  '
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.03'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.20%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +5.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.09'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  -2.75%  '
$ws.Range('E33').Value = '  +4.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '27.17'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '46.47'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.97%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '49.07'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('E41').Value = '  -8.82%  '
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('E43').Value = '  -4.42%  '
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '384.22'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('D47').Value = '2.679.31'
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '133.82'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('E51').Value = '  +1.26%  '